$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.189712882041931
$ws.Range("B1").Value = 2.365930080413818
$ws.Range("C1").Value = 4.284689903259277
$ws.Range("D1").Value = 2.889295101165771
$ws.Range("E1").Value = 1.125526905059814
